$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14, shifting existing rows 14..156 down to 15..157
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with this week's record
$ws.Range("A14").Value = 8
$ws.Range("B14").Value = "Terminal La Palmera de La Serena"
$ws.Range("C14").Value = "Coquimbo"
$ws.Range("D14").Value = 44761
$ws.Range("D14").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E14").Value = 4
$ws.Range("F14").Value = 100112040
$ws.Range("G14").Value = "Cilantro"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 3200
$ws.Range("K14").Value = 2000
$ws.Range("L14").Value = 2500
$ws.Range("M14").Value = 2250
$ws.Range("N14").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O14").Value = "Provincia del Elquí"
$ws.Range("P14").Value = 1500
$ws.Range("Q14").Value = 1.5
$ws.Range("R14").Value = "Hortaliza"
